# Add a new "TextBox 1" shape (author credit line) to slide 1, matching the
# upload diff: a borderless, no-autowrap textbox positioned top-right of the
# canvas containing "Austin Mestayter, Nick Wetta, Brian Gates".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Position/size in points (EMU / 12700): off x=8102226,y=0 ext cx=4089774,cy=369332
$left   = 8102226 / 12700
$top    = 0 / 12700
$width  = 4089774 / 12700
$height = 369332 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 1"

# <a:noFill/> on the shape
$tb.Fill.Visible = 0

# <a:bodyPr wrap="none" ...><a:spAutoFit/></a:bodyPr>
$tf = $tb.TextFrame
$tf.WordWrap = 0
$tf.AutoSize = 1

# Build the run-split text exactly like the authored runs.
$tr = $tf.TextRange
$tr.Text = "Austin "
$tr.InsertAfter("Mestayter") | Out-Null
$tr.InsertAfter(", Nick ") | Out-Null
$tr.InsertAfter("Wetta") | Out-Null
$tr.InsertAfter(", Brian Gates") | Out-Null
